# This script re-applies the "cfr_species" predictor re-run of the regression
# results table on the "Species count models" sheet, replacing the previous
# "biom_species" predictor rows (and re-numbering/refreshing the affected model
# rows) with freshly computed coefficients, per the "Re-run all analyses dated
# today" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "cfr_species"
$ws.Range("F2").Value = "cfr_species"

# Row 3
$ws.Range("B3").Value = "cfr_species"
$ws.Range("C3").Value = "hh"
$ws.Range("D3").Value = 0.3501109778881073
$ws.Range("E3").Value = 0.091486096382141113
$ws.Range("F3").Value = "cfr_species"
$ws.Range("H3").Value = 0.16713878512382507
$ws.Range("I3").Value = 0.53308320045471191

# Row 4
$ws.Range("B4").Value = "cfr_species"
$ws.Range("C4").Value = "mkt"
$ws.Range("D4").Value = 0.37285584211349487
$ws.Range("E4").Value = 0.085471294820308685
$ws.Range("F4").Value = "cfr_species"
$ws.Range("H4").Value = 0.2019132524728775
$ws.Range("I4").Value = 0.54379844665527344

# Row 5
$ws.Range("A5").Value = "consumption_species"
$ws.Range("B5").Value = "catch_species"
$ws.Range("C5").Value = "base"
$ws.Range("D5").Value = 0.89913344383239746
$ws.Range("E5").Value = 0.020791295915842056
$ws.Range("F5").Value = "catch_species"
$ws.Range("H5").Value = 0.85755085945129395
$ws.Range("I5").Value = 0.94071602821350098

# Row 6
$ws.Range("A6").Value = "consumption_species"
$ws.Range("B6").Value = "catch_species"
$ws.Range("D6").Value = 0.89727693796157837
$ws.Range("E6").Value = 0.01982947438955307
$ws.Range("F6").Value = "catch_species"
$ws.Range("H6").Value = 0.85761797428131104
$ws.Range("I6").Value = 0.9369359016418457

# Row 7
$ws.Range("A7").Value = "consumption_species"
$ws.Range("B7").Value = "catch_species"
$ws.Range("D7").Value = 0.89582264423370361
$ws.Range("E7").Value = 0.019301192834973335
$ws.Range("F7").Value = "catch_species"
$ws.Range("H7").Value = 0.85722023248672485
$ws.Range("I7").Value = 0.93442505598068237

# Row 8
$ws.Range("A8").Value = "sold_species"
$ws.Range("B8").Value = "catch_species"
$ws.Range("C8").Value = "base"
$ws.Range("D8").Value = 0.2708798348903656
$ws.Range("E8").Value = 0.039570342749357224
$ws.Range("F8").Value = "catch_species"
$ws.Range("H8").Value = 0.19173914194107056
$ws.Range("I8").Value = 0.35002052783966064

# Row 9
$ws.Range("A9").Value = "sold_species"
$ws.Range("C9").Value = "hh"
$ws.Range("D9").Value = 0.27627730369567871
$ws.Range("E9").Value = 0.036794524639844894
$ws.Range("H9").Value = 0.20268824696540833
$ws.Range("I9").Value = 0.3498663604259491

# Row 10
$ws.Range("A10").Value = "sold_species"
$ws.Range("C10").Value = "mkt"
$ws.Range("D10").Value = 0.28031918406486511
$ws.Range("E10").Value = 0.038133770227432251
$ws.Range("H10").Value = 0.20405164361000061
$ws.Range("I10").Value = 0.35658672451972961

# Row 11
$ws.Range("B11").Value = "cfr_species"
$ws.Range("C11").Value = "base"
$ws.Range("D11").Value = 0.29589375853538513
$ws.Range("E11").Value = 0.081302262842655182
$ws.Range("F11").Value = "cfr_species"
$ws.Range("G11").Value = "no_effort"
$ws.Range("H11").Value = 0.13328923285007477
$ws.Range("I11").Value = 0.45849829912185669

# Row 12
$ws.Range("B12").Value = "cfr_species"
$ws.Range("D12").Value = 0.30266293883323669
$ws.Range("E12").Value = 0.078538224101066589
$ws.Range("F12").Value = "cfr_species"
$ws.Range("G12").Value = "no_effort"
$ws.Range("H12").Value = 0.14558649063110352
$ws.Range("I12").Value = 0.45973938703536987

# Row 13
$ws.Range("B13").Value = "cfr_species"
$ws.Range("C13").Value = "mkt"
$ws.Range("D13").Value = 0.32400384545326233
$ws.Range("E13").Value = 0.071458756923675537
$ws.Range("F13").Value = "cfr_species"
$ws.Range("G13").Value = "no_effort"
$ws.Range("H13").Value = 0.18108633160591125
$ws.Range("I13").Value = 0.4669213593006134

# Row 14
$ws.Range("A14").Value = "sold_species"
$ws.Range("B14").Value = "cfr_species"
$ws.Range("C14").Value = "base"
$ws.Range("D14").Value = 0.063580483198165894
$ws.Range("E14").Value = 0.035771388560533524
$ws.Range("F14").Value = "cfr_species"
$ws.Range("G14").Value = "no_effort"
$ws.Range("H14").Value = -0.0079622939229011536
$ws.Range("I14").Value = 0.13512325286865234

# Row 15
$ws.Range("A15").Value = "sold_species"
$ws.Range("B15").Value = "cfr_species"
$ws.Range("C15").Value = "hh"
$ws.Range("D15").Value = 0.053387299180030823
$ws.Range("E15").Value = 0.035245683044195175
$ws.Range("F15").Value = "cfr_species"
$ws.Range("G15").Value = "no_effort"
$ws.Range("H15").Value = -0.017104066908359528
$ws.Range("I15").Value = 0.12387866526842117

# Row 16
$ws.Range("B16").Value = "cfr_species"
$ws.Range("C16").Value = "mkt"
$ws.Range("D16").Value = 0.0559505894780159
$ws.Range("E16").Value = 0.034688640385866165
$ws.Range("F16").Value = "cfr_species"
$ws.Range("G16").Value = "no_effort"
$ws.Range("H16").Value = -0.013426691293716431
$ws.Range("I16").Value = 0.12532787024974823

# Row 23
$ws.Range("B23").Value = "cfr_species"
$ws.Range("F23").Value = "cfr_species"

# Row 24
$ws.Range("B24").Value = "cfr_species"
$ws.Range("F24").Value = "cfr_species"

# Row 25
$ws.Range("B25").Value = "cfr_species"

# Row 26
$ws.Range("B26").Value = "cfr_species"
$ws.Range("F26").Value = "cfr_species"

# Row 27
$ws.Range("B27").Value = "cfr_species"

# Row 28
$ws.Range("B28").Value = "cfr_species"
$ws.Range("F28").Value = "cfr_species"

# Row 29
$ws.Range("B29").Value = "cfr_species"

# Row 30
$ws.Range("B30").Value = "cfr_species"
$ws.Range("F30").Value = "cfr_species"

# Row 31
$ws.Range("B31").Value = "cfr_species"
$ws.Range("F31").Value = "cfr_species"

# Row 32
$ws.Range("B32").Value = "cfr_species"

# Row 33
$ws.Range("B33").Value = "cfr_species"
$ws.Range("F33").Value = "cfr_species"

# Row 34
$ws.Range("B34").Value = "cfr_species"

# Row 35
$ws.Range("B35").Value = "cfr_species"
$ws.Range("F35").Value = "cfr_species"

# Row 36
$ws.Range("B36").Value = "cfr_species"
